# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted into the dataset at row 513,
# pushing the previously-existing rows 513-576 down to 514-577.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 513 (shifts old 513..576 down to 514..577)
$ws.Rows.Item(513).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A513").Value = 4
$ws.Range("B513").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C513").Value = "Los Lagos"
$ws.Range("D513").Value = 45142
$ws.Range("E513").Value = 10
$ws.Range("F513").Value = 100112023
$ws.Range("G513").Value = "Brócoli"
$ws.Range("H513").Value = "Sin especificar"
$ws.Range("I513").Value = "Primera"
$ws.Range("J513").Value = 1500
$ws.Range("K513").Value = 1500
$ws.Range("L513").Value = 1500
$ws.Range("M513").Value = 1500
$ws.Range("N513").Value = "$/unidad"
$ws.Range("O513").Value = "Región Metropolitana"
$ws.Range("P513").Value = 1500
$ws.Range("Q513").Value = 1
$ws.Range("R513").Value = "Hortaliza"
